$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New ShipmentTrackNum / PackageTrackNum values for rows 2-22 (column C, and
# column D where it mirrors column C), replacing the old tracking numbers.
$newValues = @(
    "320018594180",
    "320018594190",
    "320018594227",
    "320018594249",
    "320018594282",
    "320018594308",
    "320018594330",
    "320018594352",
    "320018594385",
    "320018594400",
    "320018594444",
    "320018594466",
    "320018594499",
    "320018594514",
    "320018594547",
    "320018594569",
    "320018594606",
    "320018594628",
    "320018594650",
    "320018594672",
    "320018594709"
)

# Rows where column D mirrors column C (PackageTrackNum equals ShipmentTrackNum)
$mirrorRows = @(5, 6, 7, 13, 14, 15, 16, 17)

# Helper (scratch) cell used to coerce numeric-looking strings into true text
# values without leaving behind any lingering number-format/style changes on
# the cells that are actually used by the sheet.
$helper = $ws.Cells.Item(200, 26)

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $i + 2
    $value = $newValues[$i]

    $helper.Formula = '="' + $value + '"'
    $helper.Copy()
    $ws.Cells.Item($row, 3).PasteSpecial(-4163)

    if ($mirrorRows -contains $row) {
        $ws.Cells.Item($row, 4).PasteSpecial(-4163)
    }
}

$helper.Clear()
$excel.CutCopyMode = 0
